$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 5 (the current row 4 holds the data that will move down to row 5)
$ws.Rows.Item(5).Insert()

# Copy the existing row 4 values down into the new row 5
for ($col = 1; $col -le 20; $col++) {
    $src = $ws.Cells.Item(4, $col)
    $dst = $ws.Cells.Item(5, $col)
    $dst.Value = $src.Value2
}
# Column D (4) is a date column formatted specially; match that formatting
$ws.Cells.Item(5, 4).NumberFormat = $ws.Cells.Item(4, 4).NumberFormat

# Update row 4 with its new values
$ws.Cells.Item(4, 4).Value = 45271   # D4 Fecha
$ws.Cells.Item(4, 13).Value = 60     # M4 Volumen
$ws.Cells.Item(4, 14).Value = 10000  # N4 Precio minimo
$ws.Cells.Item(4, 15).Value = 10000  # O4 Precio maximo
$ws.Cells.Item(4, 16).Value = 10000  # P4 Precio promedio ponderado
$ws.Cells.Item(4, 19).Value = 5000   # S4 Precio $/Kg
